$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 78, shifting existing rows 78+ down by one.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new leave entry.
$ws.Range("B78").Value = "SL(1-0-0)"
$ws.Range("H78").Value = 1
$ws.Range("K78").Value = (Get-Date -Year 2023 -Month 4 -Day 12 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("B79").Select()
